$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'دوز کلیوی دوپامین چقدر است؟'
$ws.Range("A3").Value = 'تشکیل آدیپوسیر با همه موارد زیر مشخص می‌شود به جز:'
$ws.Range("A4").Value = 'کدام یک از موارد زیر از ویژگی‌های کم‌خونی فانکونی نیست؟'
$ws.Range("A5").Value = 'نِفرُوتوکسیسیتی (سمیت کلیوی) یکی از عوارض جانبی کدام یک از داروهای سرکوب‌کننده سیستم ایمنی زیر است؟'
$ws.Range("A6").Value = 'همه موارد زیر به جز یکی، مستعد کننده گلوکوم با زاویه بسته هستند:'
$ws.Range("A7").Value = 'شما قصد انجام تزریق درون استخوانی را در یک کودک دارید. کدام یک از موارد زیر قبل از انجام تزریق نیاز به رد کردن ندارد؟'
$ws.Range("A8").Value = 'نسبت fevi/fvc در کدام مورد کاهش می‌یابد؟'
$ws.Range("A9").Value = 'اریتروپویتین توسط کدام بخش ترشح می‌شود؟ مارس 2009'
$ws.Range("A10").Value = 'کدام یک از موارد زیر نشان‌دهنده استفاده از پروتز ونیر کامل به عنوان نگهدارنده برای پل ثابت (fpd) است؟'
$ws.Range("A11").Value = 'همه موارد زیر در مورد مطالعات کوهورت صحیح است به جز -'
$ws.Range("A12").Value = 'دررفتگی آرنج (پولده البو) با کدام روش کاهش می‌یابد؟'
$ws.Range("A13").Value = 'اصطلاحی که برای توصیف حالتی که بیش از یک کدون برای یک اسید آمینه کد می‌کنند، استفاده می‌شود چیست؟'
$ws.Range("A14").Value = 'شکستگی آرامگاه‌سازان (undertakers fracture)-'
$ws.Range("A15").Value = 'زن جوانی با فشار خون بالا و سطح vma بیش از 14 میلی‌گرم در روز مراجعه می‌کند. این حالت با کدام یک از موارد زیر مرتبط است؟  
a) سرطان مدولاری تیروئید  
b) بیماری فون هیپل لینداو  
c) سندرم استورج وبر  
d) بیماری گریوز  
e) نوروفیبروماتوز'
$ws.Range("A16").Value = 'شایع‌ترین تظاهر حالت یوتیروئید بیمار چیست؟'
$ws.Range("A17").Value = 'کودک در چه زمانی به این مرحله رشدی می‌رسد؟'
$ws.Range("A18").Value = 'میزان آهن در تاج‌های استیل ضدزنگ 3m چقدر است؟'
$ws.Range("A19").Value = 'کدام یک از ویروس‌های هپاتیت زیر، یک ویروس dna است؟'
$ws.Range("A20").Value = 'یک مرد 20 ساله که الکل زیاد مصرف می‌کرد و در فعالیت‌های جنسی غیرقانونی درگیر بود، با کمبود خواب، تحریک‌پذیری و عدم خستگی از 3 هفته گذشته، تشخیص چیست؟'
$ws.Range("A21").Value = 'یک پسر ۳ ساله با تب، دیزوری و هماچوری آشکار مراجعه می‌کند. معاینه فیزیکی ناحیه سوپراپوبیک برجسته و کدر در پرکاشن را نشان می‌دهد. آزمایش ادرار گلبول‌های قرمز اما بدون پروتئینوری را نشان می‌دهد. کدام یک از موارد زیر محتمل‌ترین تشخیص است؟'
$ws.Range("A22").Value = 'استئوسارکوم با میوزیت اسیفیکانز در رادیولوژی توسط کدام مورد تفکیک می‌شود؟'
$ws.Range("A23").Value = 'کربوهیدرات c در استرپتوکوک همولیتیکوس برای چه چیزی مهم است؟'
$ws.Range("A24").Value = 'یک زن ۴۸ ساله که از منوراژی شدید (dub) رنج می‌برد، تحت عمل هیسترکتومی قرار گرفته است. او مایل به دریافت هورمون درمانی جایگزین است. معاینه فیزیکی و پستان طبیعی است، اما اشعه ایکس پوکی استخوان را نشان می‌دهد. درمان انتخابی چیست؟'
$ws.Range("A25").Value = 'هنگامی که یک سلول به رده‌ای متفاوت تبدیل می‌شود، این توانایی به عنوان چه شناخته می‌شود؟'
$ws.Range("A26").Value = 'همه موارد زیر از مزایای فضا نگهدارنده‌های متحرک هستند به جز'
$ws.Range("A27").Value = 'یک کودک 10 ساله پس از 2 روز اسهال دچار هماچوری می‌شود. لام خون گلبول‌های قرمز قطعه‌قطعه و ترومبوسیتوپنی را نشان می‌دهد. سونوگرافی بزرگ‌شدگی قابل توجه هر دو کلیه را نشان می‌دهد. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A28").Value = 'تومور ناشی از مخاط بویایی بینی چیست؟'
$ws.Range("A29").Value = 'گاستروشیزیس با کدام یک از موارد زیر مرتبط است؟'
$ws.Range("A30").Value = 'اتیامبوتول باعث کدام یک از موارد زیر می‌شود؟'
$ws.Range("A31").Value = 'یک خانم 56 ساله با شکایت اصلی سردرد، سرگیجه و خارش عمومی به خصوص پس از دوش آب گرم به درمانگاه مراجعه می‌کند. همچنین از سوزش شدید در دست‌ها و پاها شکایت دارد. او گزارش می‌دهد که مصرف آسپرین این مشکل را تسکین می‌دهد. معاینه فیزیکی: اسپلنومگالی با فشار خون بالا. آزمایش خون نشان می‌دهد: هموگلوبین 20.1 گرم در دسی‌لیتر، هماتوکریت 60%، گلبول سفید 15800، شمارش پلاکت 500000، اریتروپویتین پایین، اشباع اکسیژن 98%، آلکالین فسفاتاز لوکوسیتی افزایش یافته. بیمار با جهش jak2v617f تشخیص داده شده است. کدام یک از اسلایدهای پاتولوژی زیر با این شرایط بالینی مطابقت دارد؟'
$ws.Range("A32").Value = 'لایه های شبکه ای عمیق درم چه نوع بافت پیوندی هستند؟'
$ws.Range("A33").Value = 'کدام یک از داروهای خط اول ضد سل در بارداری منع مصرف دارد؟'
$ws.Range("A34").Value = 'زخم شدن فرج معمولاً در همه موارد زیر دیده می‌شود به جز:'
$ws.Range("A35").Value = 'همه موارد زیر در مورد توهمات الکلی صحیح هستند به جز'
$ws.Range("A36").Value = 'کدام اندام از اسیدهای چرب استفاده نمی‌کند؟'
$ws.Range("A37").Value = 'همه موارد زیر جزو men-1 هستند به جز'
$ws.Range("A38").Value = 'اصلی ترین سورفکتانت کدام است؟'
$ws.Range("A39").Value = 'در مورد تحلیل پیوند در اختلالات ژنتیکی خانوادگی چه صحیح است؟'
$ws.Range("A40").Value = 'طبق برنامه ایمن‌سازی، واکسن هپاتیت a در چه سنی توصیه می‌شود؟'
$ws.Range("A41").Value = 'مجموعه‌های کلاسپ که حرکت عملکردی پروتز را تطبیق می‌دهند، برای حل کدام نگرانی طراحی شده‌اند؟'
$ws.Range("A42").Value = 'آنستوموز لنفاوی-وریدی برای چه شرایطی انجام می‌شود؟'
$ws.Range("A43").Value = 'کاهش کدام یک از موارد زیر محتمل‌ترین علت ادم محیطی در یک بیمار با سابقه طولانی‌مدت اعتیاد به الکل و بیماری کبدی است؟'
$ws.Range("A44").Value = 'بیشترین مقدار ویتامین c در کدام گزینه یافت می‌شود؟'
$ws.Range("A45").Value = 'مخرج در نرخ مرگ و میر خاص ناشی از بیماری عروق کرونر چیست؟'
$ws.Range("A46").Value = 'کدام یک از موارد زیر از دستیاران دندانپزشکی غیرعملیاتی محسوب می‌شود؟'
$ws.Range("A47").Value = 'عفونت دراکونکولوزیس چگونه اتفاق می‌افتد؟'
$ws.Range("A48").Value = 'در مورد سندرم درسلر همه موارد زیر صحیح است، به جز -'
$ws.Range("A49").Value = 'ناحیه بویایی در بینی کجاست؟'
$ws.Range("A50").Value = 'یک بیمار 55 ساله با سابقه انفارکتوس میوکارد، در 6 ماه گذشته از آسپرین 75 میلی‌گرم استفاده می‌کند. او نیاز به کشیدن دندان پوسیده دارد. دندانپزشک باید چه کاری انجام دهد؟'
$ws.Range("A51").Value = 'بالاترین میزان بروز کم خونی در مناطق گرمسیری به چه علتی است؟'
$ws.Range("A52").Value = 'برون‌ده قلبی در تمام شرایط زیر کاهش می‌یابد به جز'
$ws.Range("A53").Value = 'کدام یک از داروهای زیر در کاهش بازسازی میوکارد در نارسایی احتقانی قلب مؤثر نیست؟'
$ws.Range("A54").Value = 'کدام یک از داروهای زیر باعث پریکاردیت نمی‌شود؟'
$ws.Range("A55").Value = 'اولین نشانه افزایش فشار داخل جمجمه چیست؟'
$ws.Range("A56").Value = 'سندرم شوک سمی استرپتوکوکی ناشی از چیست؟'
$ws.Range("A57").Value = 'پانانسفالیت اسکلروزان تحت حاد عارضه کدام یک از موارد زیر است؟'
$ws.Range("A58").Value = 'در یک بیمار با po2 برابر 85 mmhg، pco2 برابر 50 mmhg، ph برابر 7.2 و hco3 برابر 32 meq/1، بیمار از چه شرایطی رنج می‌برد؟'
$ws.Range("A59").Value = 'اریتم پالمار در کدام مورد مشاهده می‌شود؟'
$ws.Range("A60").Value = 'آزمایش فلورانس برای تشخیص چه چیزی استفاده می‌شود؟'
$ws.Range("A61").Value = 'کدام یک از موارد زیر یک بیماری وابسته به y است؟'
$ws.Range("A62").Value = 'کدام یک از موارد زیر در مورد آنتاگونیسم رقابتی برگشت‌پذیر صحیح است؟'
$ws.Range("A63").Value = 'در مورد پلی‌مورفیسم کدام گزینه صحیح است؟'
$ws.Range("A64").Value = 'در پاپیلومای ویلیوس رکتوم کدام یک از دست می‌رود؟'
$ws.Range("A65").Value = 'بیمار با هایپربیلی روبینمی غیرکونژوگه و سطح بالای اوروبیلینوژن در ادرار مراجعه می‌کند. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A66").Value = 'بهترین شاخص برای تعیین سطح کلی آلودگی هوا در یک منطقه شهری کدام است؟'
$ws.Range("A67").Value = 'از چند سانتی‌متر باز شدن دهانه رحم، پارتوگرام به صورت منظم ترسیم می‌شود؟'
$ws.Range("A68").Value = 'شایع‌ترین نقص دیواره بین دهلیزی در سندرم داون کدام است؟'
$ws.Range("A69").Value = 'استروئیدهای موضعی در درمان کدام مورد استفاده می‌شوند؟'
$ws.Range("A70").Value = 'امیل دورکیم با کار بر روی کدام یک از شرایط زیر در روانپزشکی مرتبط است؟'
$ws.Range("A71").Value = 'جسم بار در کدام مورد مشاهده می‌شود؟'
$ws.Range("A72").Value = 'کدام ناهنجاری کروموزومی با سندرم چشم گربه مرتبط است؟'
$ws.Range("A73").Value = 'یک کودک 2 ساله پس از نوشیدن مایع ضدعفونی کننده به اورژانس آورده شد. در معاینه، او در حالت کما بود، مردمک‌های چشم تنگ شده بود، نبض 120 در دقیقه؛ تعداد تنفس: 40 در دقیقه، ادرار کم بود و رنگ آن در تماس با هوا به سبز تغییر کرد. تشخیص چیست؟'
$ws.Range("A74").Value = 'در دوز پایین، افزایش خطی در غلظت پلاسمایی دارو و در دوز بالاتر، افزایش شدید در غلظت دارو، چه نوع سینتیکی است؟'
$ws.Range("A75").Value = 'پرتو درمانی داخل حین عمل به کدام عضو داده می‌شود؟'
$ws.Range("A76").Value = 'همه قندهای احیا کننده هستند به جز'
$ws.Range("A77").Value = 'هرس شریان‌های ریوی در کدام قسمت مشاهده می‌شود؟'
$ws.Range("A78").Value = 'وقتی اطلاعاتی که بعداً به خاطر سپرده می‌شوند توسط اطلاعاتی که قبلاً یاد گرفته شده‌اند مختل می‌شوند، به آن چه می‌گویند؟'
$ws.Range("A79").Value = 'یافته پاتوگنومونیک بیماری هاری کدام است؟ مارس 2012'
$ws.Range("A80").Value = 'کدام یک از علائم زیر با turp بهبود نمی یابد؟'
$ws.Range("A81").Value = 'داده‌های کیفی را می‌توان با کدام روش نمایش داد؟'
$ws.Range("A82").Value = 'کدام یک از موارد زیر دیستروفی کمربند اندامی نیست؟'
$ws.Range("A83").Value = 'درمان جایگزینی هورمونی (hrt) چه چیزی را بهبود می‌بخشد؟'
$ws.Range("A84").Value = 'کدام یک از موارد زیر تحت تأثیر بیماری‌های سیستمیک قرار نمی‌گیرد؟'
$ws.Range("A85").Value = 'کیست زیرملتحمه در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A86").Value = 'همه موارد زیر آگونیست‌های دوپامینرژیک مورد استفاده در پارکینسونیسم هستند به جز'
$ws.Range("A87").Value = 'کدام یک از موارد زیر باعث دورسی فلکشن پا می‌شود؟'
$ws.Range("A88").Value = 'کدام یک از موارد زیر باعث پاراپلژی اسپاستیک گرمسیری می‌شود؟'
$ws.Range("A89").Value = 'موارد نشان‌دهنده زایمان فوری نوزاد دوم در دوقلوها شامل همه موارد زیر به جز کدام است؟'
$ws.Range("A90").Value = 'کدام یک از موارد زیر نشانه هیداتیدوز ریوی نیست؟'
$ws.Range("A91").Value = 'کدام گزینه اثر وابسته به زمان در کشتن باکتری را نشان می‌دهد؟'
$ws.Range("A92").Value = 'تمامی عوامل زیر خطر سمیت کلیوی آمینوگلیکوزیدها را افزایش می‌دهند به جز -'
$ws.Range("A93").Value = 'همه موارد زیر در مورد سوفل فیزیولوژیک صحیح است به جز'
$ws.Range("A94").Value = 'هیپوناترمی در تمام موارد زیر دیده می‌شود به جز:'
$ws.Range("A95").Value = 'رها کردن (یا فراموش کردن) یک ابزار یا اسفنج در شکم بیمار در حین جراحی و بستن عمل جراحی، چه نوع خطایی محسوب می‌شود؟'
$ws.Range("A96").Value = 'عارضه زودرس tipss کدام است؟'
$ws.Range("A97").Value = 'اجزای خون زیر را به ترتیب نزولی مدت زمان ماندگاری آنها مرتب کنید:  
a. پلاکت در دمای 20 درجه سانتی‌گراد  
b. ffp در دمای 18- درجه سانتی‌گراد  
c. خون کامل با محلول cpda-1 در دمای 4 درجه سانتی‌گراد  
d. خون کامل با محلول cpd در دمای 4 درجه سانتی‌گراد'
$ws.Range("A98").Value = 'واکسن dpt منجمد شده باید -'
$ws.Range("A99").Value = 'کدام یک از ویژگی‌های ترک نیکوتین نیست؟'
$ws.Range("A100").Value = 'کدام یک از موارد زیر اولین و پایدارترین علامت تومور گلوموس است؟'
$ws.Range("A101").Value = 'فعال شدن رنین باعث تحریک چه چیزی می‌شود؟'
$ws.Range("A102").Value = 'دیسپلازی کامل حلزون گوش به چه نامی شناخته می‌شود؟'
$ws.Range("A103").Value = 'هیپوپلازی ریوی همراه با مشکلات ادراری با کدام یک از موارد زیر مرتبط است؟'
$ws.Range("A104").Value = 'روش انتخابی برای بررسی عمق نفوذ و مرحله گره‌های لنفاوی در کارسینوم رکتوم کدام است؟'
$ws.Range("A105").Value = 'آنزیم محدودکننده سرعت در سنتز اسیدهای صفراوی کدام است؟'
$ws.Range("A106").Value = 'در یک آزمایش غربالگری برای سرطان پستان، حساسیت 90٪ و ویژگی 98٪ است. احتمال اینکه یک بیمار مبتلا به سرطان پستان در غربالگری انجام شده در دو سال متوالی نتیجه منفی بگیرد چقدر است؟'
$ws.Range("A107").Value = 'یک کارگر کارخانه 36 ساله دچار پلاک‌های خارش‌دار حلقوی و پوسته‌دار در هر دو کشاله ران شده است. استفاده از پماد کورتیکواستروئید باعث تسکین موقت شد اما پلاک‌ها همچنان در محیط گسترش یافتند. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A108").Value = 'خط نلتون خطی است که کدام نقاط را به هم وصل می‌کند؟'
$ws.Range("A109").Value = 'کدام بیماری تحت نظارت بین‌المللی است؟'
$ws.Range("A110").Value = 'تمام موارد زیر در مورد رانولا صحیح هستند به جز'
$ws.Range("A111").Value = 'hmb 45 یک نشانگر ایمونوهیستولوژیک برای کدام مورد است؟'
$ws.Range("A112").Value = 'شکست روش پیشگیری از بارداری 10/100 زن-سال به چه معناست؟'
$ws.Range("A113").Value = 'الیاف بالارونده قشر مخچه کدامند؟'
$ws.Range("A114").Value = 'یک پسر 8 ساله با وزن 30 کیلوگرم تحت عمل برداشت تومور ویلمز قرار دارد. هموگلوبین او 12 گرم در دسی لیتر است. اگر آستانه تزریق خون 8 گرم در دسی لیتر باشد، میزان خونریزی مجاز چقدر است؟'
$ws.Range("A115").Value = 'پروفایل بیوفیزیکال جنین که برای ارزیابی سلامت جنین استفاده می‌شود، شامل کدام مورد زیر نمی‌شود؟'
$ws.Range("A116").Value = 'رینوفیما یک تومور با رشد آهسته است که دارای پاتولوژی زیر می‌باشد:'
$ws.Range("A117").Value = 'در مورد ایمونوگلوبولین‌ها کدام گزینه صحیح است؟'
$ws.Range("A118").Value = 'کدام یک از گزینه های زیر داروی انتخابی برای درمان سندرم ترشح نامناسب هورمون ضدادراری است؟'
$ws.Range("A119").Value = 'در مورد فئوکروموسیتوما، همه موارد زیر صحیح هستند به جز -'
$ws.Range("A120").Value = 'برای نوزادی که به طور خودبه‌خودی نفس می‌کشد، از کدام مدار تنفسی باید استفاده کرد؟'
$ws.Range("A121").Value = 'روش ترجیحی شانت برای بیمار مبتلا به فشار خون پورت با ریسک جراحی قابل قبول و عملکرد کافی کبد کدام است؟'
$ws.Range("A122").Value = 'تظاهرات بالینی شامل گریه خشن، فتق نافی، هیپوتونی، لکه‌لکه شدن پوست، لتارژی و زردی طولانی‌مدت در کدام بیماری دیده می‌شود؟'
$ws.Range("A123").Value = 'اسید وانیلی ماندلیک (vma) در ادرار در کدام مورد افزایش می یابد؟'
$ws.Range("A124").Value = 'موارد زیر از اجزای بیماری کاوازاکی هستند به جز'
$ws.Range("A125").Value = 'نادرست در مورد بیماری شواخمن'
$ws.Range("A126").Value = 'نمونه‌گیری طبقه‌بندی شده برای چه نوع داده‌ای ایده‌آل است؟'
$ws.Range("A127").Value = 'شایع‌ترین محل کیست برانشیال در کجاست؟'
$ws.Range("A128").Value = 'مقدار کمی آتروپین به دیفنوکسیلات اضافه می‌شود زیرا'
$ws.Range("A129").Value = 'یک زن 40 ساله با تنگی نفس در هنگام فعالیت و تپش قلب مراجعه کرده است. نوار قلب نشان دهنده فیبریلاسیون دهلیزی با ضربان بطنی سریع است. اکوکاردیوگرافی تنگی شدید دریچه میترال با لخته در زائده دهلیز چپ را نشان می‌دهد. کدام یک از موارد زیر توصیه نمی‌شود؟'
$ws.Range("A130").Value = 'دمانس‌های فرونتوتمپورال شامل کدام مورد نمی‌شوند؟'
$ws.Range("A131").Value = 'اگر حد اطمینان افزایش یابد، آنگاه:'
$ws.Range("A132").Value = 'بیمار دو روز پس از یک دوره مالاریا دچار قرمزی چشم شده است. علت احتمالی چیست؟'
$ws.Range("A133").Value = 'دی‌اکسید کربن در خون بیشتر به چه صورت حمل می‌شود؟'
$ws.Range("A134").Value = 'درگیری پوستی در نکروز اپیدرمال سمی چقدر است؟'
$ws.Range("A135").Value = 'در درجه‌بندی تراخم، التهاب تراخمایی فولیکولی به عنوان وجود چه چیزی تعریف می‌شود؟'
$ws.Range("A136").Value = 'یک کودک ۵ ساله از خوابیدن در تخت خود امتناع می‌کند و ادعا می‌کند که هیولاهایی در کمد او وجود دارند و کابوس می‌بیند. والدین به او اجازه می‌دهند در تخت آنها بخوابد تا از جیغ‌زدن های اجتناب‌ناپذیر او جلوگیری کنند. والدین اشاره می‌کنند که کودک به‌خوبی می‌خوابد و فقط با طلوع آفتاب بیدار می‌شود. اختلال خوابی که بیشترین تطابق را با این شرح حال دارد، انتخاب کنید.'
$ws.Range("A137").Value = 'در مورد سینوس کرونری کدام گزینه صحیح است؟'
$ws.Range("A138").Value = 'یک مرد 40 ساله مبتلا به فشار خون بالا با سردرد ناگهانی و تغییر وضعیت هوشیاری در بیمارستان بستری شد. در معاینه، فشار خون او 220/110 میلی‌متر جیوه مشاهده شد و بیمار چهار ساعت بعد فوت کرد. یافته پاتولوژیک احتمالی در کلیه‌های او چیست؟'
$ws.Range("A139").Value = 'کدام یک از موارد زیر، منع مصرف برای واکسن opv محسوب نمی‌شود؟'
$ws.Range("A140").Value = 'یک فرد الکلی 35 ساله با اسهال و بثورات پوسته‌دار در گردن، دست‌ها و پاها مراجعه می‌کند. کدام یک از موارد زیر باعث بهبودی او می‌شود؟'
$ws.Range("A141").Value = 'شایع‌ترین علت کراتیت در کاربران لنزهای تماسی نرم چیست؟'
$ws.Range("A142").Value = 'یک دانشجوی پزشکی در معرض یک بیمار مبتلا به سل قرار گرفته و تست توبرکولین (ppd) مثبت داشته است، اما عکس قفسه سینه طبیعی نشان داده است. او تحت یک دوره 6 ماهه درمان پیشگیرانه قرار می‌گیرد، اما بعداً دچار نوروپاتی محیطی می‌شود. کدام یک از ویتامین‌های زیر برای درمان نوروتوکسیسیتی در نظر گرفته می‌شود؟'
$ws.Range("A143").Value = 'اسیدهای صفراوی از چه چیزی تشکیل می‌شوند؟'
$ws.Range("A144").Value = 'کدام ماده بیهوشی بیشترین mac را دارد؟'
$ws.Range("A145").Value = 'در یک حیوان فاقد مغز استخوان، سلول‌های سری میلوئیدی تزریق می‌شوند. پس از دوره انکوباسیون کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Range("A146").Value = 'مثال‌هایی از حساسیت نوع ۴ کدام است؟'
$ws.Range("A147").Value = 'اسمیر محیطی در اسفروسیتوز ارثی، اسفروسیت‌ها را نشان می‌دهد که -'
$ws.Range("A148").Value = 'درصد dna میتوکندریایی از کل dna سلولی چقدر است؟'
$ws.Range("A149").Value = 'یک زن ۵۴ ساله با پاپول‌های اریتماتوز پراکنده در سمت راست پیشانی مراجعه می‌کند. او در دو روز گذشته درد "سوزشی" و حساسیت بیش از حد در آن ناحیه داشته است. درد او به ناحیه سوپرااوربیتال راست و قسمت راست پشت بینی محدود می‌شود. او سردرد، تغییرات وضعیت روانی یا عفونت‌های اخیر را انکار می‌کند. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A150").Value = 'کدام یک از علائم زیر را در یک بیمار مبتلا به پیلونفریت حاد انتظار دارید؟'
$ws.Range("A151").Value = 'در مورد اختلال اسید-باز، همه موارد صحیح هستند به جز:'
$ws.Range("A152").Value = 'کدام یک از موارد زیر واکسن زنده نیست؟'
$ws.Range("A153").Value = 'اجسام شیلر-دووال در کدام یک دیده می‌شوند؟'
$ws.Range("A154").Value = 'مقدار ریشه عصبی رفلکس زانو چیست؟'
$ws.Range("A155").Value = 'برآمدگی سرویکال. کدام گزینه صحیح است؟'
$ws.Range("A156").Value = 'یک زن ۲۸ ساله از درد قفسه سینه به مدت ۱.۵ ماه شکایت دارد. سه ماه قبل از مراجعه فعلی، او دچار سرفه همراه با خلط شده بود که در صبح زود بدتر می‌شد. تنگی نفس یا هموپتیزی وجود نداشت. حدود دو ماه پس از شروع سرفه، او شروع به درد متناوب پلورتیک (بدتر با سرفه یا تنفس عمیق) در سمت چپ قفسه سینه، افزایش خستگی و کاهش وزن ۳ کیلوگرمی کرد اما از دست دادن اشتها نداشت. در سمع: کراکل‌های ظریف در نواحی سوپرااسکاپولار دو طرفه وجود داشت (چپ بیشتر از راست). لوب فوقانی چپ در پرکاشن کدر بود. همه موارد زیر جزء محتویات محیط کلاسیک مورد استفاده برای بیماری فوق هستند به جز:'
$ws.Range("A157").Value = 'کدام یک از موارد زیر باقیماندهٔ نوتوکورد است؟'
$ws.Range("A158").Value = 'وزیکول‌های ماتریکس نقش مهمی در معدنی‌سازی چه چیزی دارند؟'
$ws.Range("A159").Value = 'در درمان کودکان مبتلا به عقب‌ماندگی ذهنی، معمولاً مشاهده می‌شود که آن‌ها:'
$ws.Range("A160").Value = 'نالتروکسان در وابستگی به مواد افیونی برای چه موردی استفاده می‌شود؟'
$ws.Range("A161").Value = 'کدام یک از ارگانیسم‌های زیر نمی‌تواند باعث عفونت جدی در بیماری گرانولوماتوز مزمن شود؟'
$ws.Range("A162").Value = 'وایت ویتریول چیست؟'
$ws.Range("A163").Value = 'درصد آب بدن در کدام گروه بیشتر است؟'
$ws.Range("A164").Value = 'شکستگی مونتگیا چیست؟ (d. repeat 2012)'
$ws.Range("A165").Value = 'مکانیسم اصلی هایپرپیرکسی ناشی از آتروپین شامل موارد زیر است:'
$ws.Range("A166").Value = 'ایمنی غیرفعال مادرزادی در کدام یک از موارد زیر یافت نمی‌شود؟'
$ws.Range("A167").Value = 'مانور لوست در زایمان چه موردی استفاده می‌شود؟'
$ws.Range("A168").Value = 'بیمار بالغ با لکه‌های پوسته‌دار کم‌رنگ روی قفسه سینه و کمر مراجعه می‌کند؛ بررسی لامپ وود از ضایعات چه رنگی را نشان می‌دهد؟'
$ws.Range("A169").Value = 'یک انفارکتوس سفید به چه دلیل رخ می‌دهد؟'
$ws.Range("A170").Value = 'فرضیه دو مرحله ای نادسون برای کدام مورد است؟'
$ws.Range("A171").Value = 'هالوفانترین برای چه موردی استفاده می‌شود؟'
$ws.Range("A172").Value = 'کدام یک از موارد زیر مربوط به حالب چپ نیست؟'
$ws.Range("A173").Value = 'کدام دارو باعث کاهش فشار داخل چشم و آپنه در نوزادان می‌شود؟'
$ws.Range("A174").Value = 'یک مرکز بهداشت اولیه (phc) در هفته 40 تا 50 مورد بیماری در جامعه گزارش می‌کند. این هفته بر اساس سوابق قبلی، 48 مورد گزارش شده است که وضعیت عادی است. این وضعیت چه نامیده می‌شود؟'
$ws.Range("A175").Value = 'بافت لنفوئید در کدام تومور پاروتید دیده می‌شود؟'
$ws.Range("A176").Value = 'درمان آلکاپتونوریا:'
$ws.Range("A177").Value = 'کدام عامل ضد دیابت در هر دو نوع دیابت نوع 1 و نوع 2 استفاده می‌شود؟'
$ws.Range("A178").Value = 'کدام یک از موارد زیر از ویژگی‌های معمول آسکاریازیس نیست؟'
$ws.Range("A179").Value = 'بزرگ‌ترین مزیت بازسازی پالپ در اندودانتیکس چیست؟'
$ws.Range("A180").Value = 'رفلکس انتروگاستریک توسط همه موارد زیر تحریک می‌شود به جز:'
$ws.Range("A181").Value = 'علامت جابجایی مایع در کدام مورد دیده می‌شود؟'
$ws.Range("A182").Value = 'تغییر چربی میکرووزیکولار در هپاتوسیت‌ها به دنبال عفونت با کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A183").Value = 'کدام یک از موارد زیر به عنوان انسولین بدون پیک شناخته می‌شود؟'
$ws.Range("A184").Value = 'جراحی استفان فاولر برای چه شرایطی انجام می‌شود؟'
$ws.Range("A185").Value = 'کدام بیماری متابولیک با انسداد روده همراه است؟'
$ws.Range("A186").Value = 'علامت آمسلر؟'
$ws.Range("A187").Value = 'تولید سورفکتانت در ریه‌ها از چه زمانی شروع می‌شود؟'
$ws.Range("A188").Value = 'نادرست در مورد لنفوم در بیماران ایدز'
$ws.Range("A189").Value = 'سریع‌ترین داروی ضد تیروئید کدام است؟'
$ws.Range("A190").Value = 'لانه‌گزینی در چه مرحله‌ای از زیگوت رخ می‌دهد؟'
$ws.Range("A191").Value = 'به طور سنتی، سرطان کیسه بیضه با کدام یک از مشاغل زیر مرتبط است؟'
$ws.Range("A192").Value = 'بلاکرهای بتا در همه موارد زیر به جز کدام مورد استفاده می‌شوند؟'
$ws.Range("A193").Value = 'آزمایش گتلر برای چیست؟'
$ws.Range("A194").Value = 'خط سیمیان در کدام مورد دیده نمی‌شود؟'
$ws.Range("A195").Value = 'شریان عبوری در کانال گویون کدام است؟'
$ws.Range("A196").Value = 'مانور والسالوا صدای سوفل را در کدام یک افزایش می‌دهد؟'
$ws.Range("A197").Value = 'کدام یک از روغن های پخت و پز زیر دارای بالاترین درصد اسید لینولئیک است؟'
$ws.Range("A198").Value = 'مری حلزونی در کدام یک از شرایط زیر مشاهده می‌شود؟'
$ws.Range("A199").Value = 'پانچایاتی راج شامل همه موارد زیر است، به جز -'
$ws.Range("A200").Value = 'یک فرد مصدوم با فلج عضله ادکتور پولیسیس، ضعف عضلات هیپوتنار و از دست دادن حس در سطح کف دستی و پشتی انگشتان چهارم و پنجم ارائه شده است. محل ضایعه کجاست؟'
$ws.Range("A201").Value = 'ساختار مارپیچ سه‌گانه در کدام ماده یافت می‌شود؟'
$ws.Range("A202").Value = 'ضایعه مغزی حلقوی تقویت‌شده در همه موارد زیر دیده می‌شود به جز'
$ws.Range("A203").Value = 'منحنی داده شده را بررسی کنید و در مورد یافته نظر دهید'
$ws.Range("A204").Value = 'تمامی عبارات زیر در مورد ویژگی‌های بافت‌شناسی آمفیزم صحیح هستند، به جز:'
$ws.Range("A205").Value = 'یک متخصص دندانپزشکی ممکن است از نظر جنایی مسئول باشد اگر مرتکب اشتباهی شود ________.'
$ws.Range("A206").Value = 'در مورد رینواسکلروما همه موارد زیر صحیح است، به جز:'
$ws.Range("A207").Value = 'بهترین گزینه درمانی برای مولوسکوم کونتاژیوزوم چیست؟'
$ws.Range("A208").Value = 'آسیب به کدام یک از ساختارهای زیر منجر به سندرم کلوور-بوسی می‌شود:'
$ws.Range("A209").Value = 'تجزیه و تحلیل گازهای خونی شریانی در ویالی حاوی هپارین باعث کاهش کدام یک از مقادیر زیر می‌شود؟'
$ws.Range("A210").Value = 'در طی تشکیل آب مروارید، خطوط عدسی ابتدا در کدام ناحیه ظاهر می‌شوند؟'
$ws.Range("A211").Value = 'تمام موارد زیر در مورد تریکینلا اسپیرالیس صحیح است به جز'
$ws.Range("A212").Value = 'تمامی عبارات زیر در مورد سرخجه مادرزادی صحیح هستند به جز'
$ws.Range("A213").Value = 'بیمار وارفارین مصرف می‌کرد و inr به 8 رسیده است. قدم بعدی قطع دارو و:'
$ws.Range("A214").Value = 'روش ترجیحی برای تأیید قرارگیری صحیح دستگاه داخل رحمی (iud) در مواردی که رشته‌های آن دیده نمی‌شود کدام است؟'
$ws.Range("A215").Value = 'در یک ترومبوس، خطوط تیره زان به دلیل چیست؟'
$ws.Range("A216").Value = 'همه موارد زیر از ویژگی‌های رادیولوژیک تنگی میترال هستند به جز –'
$ws.Range("A217").Value = 'یک بیمار مرد 28 ساله با شکایت اصلی تب، درد حاد شکم و لکوسیتوز که نشاندهنده آپاندیسیت حاد است، به بخش اورژانس منتقل میشود. لاپاراتومی اورژانس برنامهریزی شده است. با این حال، در حین جراحی یافته زیر مشاهده میشود. تعداد زیادی غدد لنفاوی مزانتریک بزرگشده وجود دارد. تشخیص احتمالی در صورت عفونی بودن میتواند باشد'
$ws.Range("A218").Value = 'پان‌سیتوپنی با مغز استخوان سلولی در همه موارد زیر دیده می‌شود به جز:'
$ws.Range("A219").Value = 'پل های بافتی در کدام یک از موارد زیر دیده می شوند؟ dnb 10'
$ws.Range("A220").Value = 'یک مرد 73 ساله با انفارکتوس تحتانی و ارتفاع st در لیدهای پرهکوردیال سمت راست، دچار افت فشار خون (فشار خون 90/70 میلی‌متر جیوه) و تاکیکاردی شده است. فشار ورید ژوگولار 10 سانتی‌متر است، صدای قلب طبیعی است، ریه‌ها تمیز هستند و اندام‌ها سرد هستند. مانیتورینگ همودینامیک مرکزی چه چیزی را نشان می‌دهد؟ برای این بیمار، پارامترهای همودینامیکی که به احتمال زیاد اعمال می‌شوند را انتخاب کنید.'
$ws.Range("A221").Value = 'بیمار به دلیل افزایش فشار داخل جمجمه‌ای ناشی از متاستازهای ثانویه در بخش مراقبت‌های ویژه بستری شده است. او دچار فلج عصب ششم سمت چپ شده است. این حالت منجر به چه می‌شود؟'
$ws.Range("A222").Value = 'فرضیه دو ضربه ای نادسن به طور مشخص با کدام یک مرتبط است؟'
$ws.Range("A223").Value = 'کدام یک از داروهای شیمی درمانی پلیمریزاسیون میکروتوبول ها را مهار می کند اما با ایجاد سرکوب مغز استخوان همراه نیست؟'
$ws.Range("A224").Value = 'عامل اتیولوژیک اولیه در ایجاد نقص فورکیشن چیست؟'
$ws.Range("A225").Value = 'یک پسر 14 ساله با سابقه خونریزی مکرر بینی مراجعه می‌کند. هموگلوبین او 6.4 گرم در دسی‌لیتر است و اسمیر محیطی نشان‌دهنده کم‌خونی نورموسیتیک هیپوکروم است. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A226").Value = 'یک کودک طبیعی کنترل کامل گردن را تا چه زمانی به دست می‌آورد؟'
$ws.Range("A227").Value = 'کیست‌های برونکوژنیک خارج‌لوب ممکن است با موارد زیر ارتباط داشته باشند به جز -'
$ws.Range("A228").Value = 'کدام یک از موارد زیر به احتمال زیاد باعث ایجاد این نوع ضایعه در بیمار مبتلا به hiv می‌شود؟'
$ws.Range("A229").Value = 'یک زن ۲۱ ساله با سابقه پتوز خفیف دو طرفه، ضعف عضلانی پروگزیمال و خستگی پذیری آسان مراجعه می‌کند. کدام یک از موارد زیر برای تشخیص این وضعیت بهترین است؟'
$ws.Range("A230").Value = 'انسداد آترواسکلروتیک کدام شریان منجر به پرفیوژن ناکافی مثانه می‌شود؟'
$ws.Range("A231").Value = 'قند خون ناشتا در یک زن باردار دیابتی باید در چه محدوده‌ای حفظ شود؟'
$ws.Range("A232").Value = 'علامت فرومنت ناشی از چیست؟'
$ws.Range("A233").Value = 'کدام بیماری از طریق تمام اجزای خون منتقل می‌شود؟'
$ws.Range("A234").Value = 'طول عمر نوتروفیل -'
$ws.Range("A235").Value = 'کدام یک نشانه خطرناک در مورد آسیب شدید سر محسوب می‌شود؟'
$ws.Range("A236").Value = 'کدام دارو باعث آگرانولوسیتوز می‌شود؟'
$ws.Range("A237").Value = 'در سندرم بدخیم نورولپتیک، علت مرگ چیست؟'
$ws.Range("A238").Value = 'مانیتورهای مدرن برای اندازه‌گیری etco2 از چه فناوری استفاده می‌کنند؟'
$ws.Range("A239").Value = 'کدام یک از شل‌کننده‌های عضلانی اسکلتی زیر باعث درد در هنگام تزریق می‌شود؟'
$ws.Range("A240").Value = 'دیواره سلولی باکتری‌های گرم مثبت حاوی چه چیزی است؟'
$ws.Range("A241").Value = 'درگیری کدام ساختار در سندرم چرگ اشتراوس غیرمعمول است؟'
$ws.Range("A242").Value = 'یک بیماری که در دوران بلوغ و در غیاب مقادیر زیاد پلاک سوپراژینجیوال رخ می‌دهد و منجر به نقص استخوانی زاویه‌ای در ناحیه مولر اول و دندان‌های پیشین می‌شود، چه نام دارد؟'
$ws.Range("A243").Value = 'نادرست در مورد سندرم مرگ ناگهانی نوزاد'
$ws.Range("A244").Value = 'تجویز محلول گلوکز در تمامی موارد زیر به جز کدام یک توصیه می‌شود؟'
$ws.Range("A245").Value = 'پوینت‌های کاغذی چگونه استریل می‌شوند؟'
$ws.Range("A246").Value = 'کدام یک از عفونت های مادرزادی زیر با حداقل خطر تراتوژنیک برای جنین همراه است؟'
$ws.Range("A247").Value = 'شایع ترین عصب آسیب دیده در ترمیم باز فتق -'
$ws.Range("A248").Value = 'در طول ورزش، افزایش تحویل اکسیژن به عضلات به دلایل زیر است به جز:'
$ws.Range("A249").Value = 'عبارت صحیح درباره آپاندیس - الف) فاقد مزانتر است ب) دارای تنیا کولی است ج) از میانه روده توسعه می یابد د) توسط شاخه آپاندیسیلار ایلئوکولیک خونرسانی می شود'
$ws.Range("A250").Value = 'کدام یک از موارد زیر نمونه‌ای از محدودیت ناتوانی است؟'
$ws.Range("A251").Value = 'کدام یک از آلودگی‌های زیر منجر به سوءجذب می‌شود؟'
$ws.Range("A252").Value = 'مادر متوجه ترشحات موکوپورولنت از نقطه سمت راست کودک 3 ماهه خود شده است. سابقه آبریزش از چشم راست کودک وجود دارد. درمان مناسب چیست؟'
$ws.Range("A253").Value = 'یک نوزاد هنگام مکیدن شیر می‌تواند بدون مشکل تنفس کند به دلیل'
$ws.Range("A254").Value = 'همه موارد زیر از عوامل بیماری‌زا مرتبط با پنومونی در الکلیسم مزمن هستند، به جز:'
$ws.Range("A255").Value = 'همه موارد زیر نمونه‌هایی از ارتقای سلامت هستند به جز'
$ws.Range("A256").Value = 'سندرم گلدن هار؟'
$ws.Range("A257").Value = '5 روز پس از جراحی cbd نشتی کوچکی مشاهده می‌شود. بهترین درمان چیست؟'
$ws.Range("A258").Value = 'داروی انتخابی برای شوک آنافیلاکتیک چیست؟'
$ws.Range("A259").Value = 'یک مرد 29 ساله یک سال پس از پیوند کلیه تحت پیگیری قرار می‌گیرد. او از نظر بالینی وضعیت خوبی دارد و پیوند به خوبی عملکرد دارد، بدون هیچ‌گونه اپیزود رد پیوند اخیر یا تغییراتی در داروهای سرکوب‌کننده سیستم ایمنی. کدام یک از عوارض زیر پیوند، محتمل‌ترین علت مرگ در این بیمار است؟'
$ws.Range("A260").Value = 'بهترین نمای رادیولوژیک برای تشخیص شکستگی اسکافوئید کدام است؟'
$ws.Range("A261").Value = 'هدف برنامه ملی کنترل نابینایی (npcb) کاهش شیوع نابینایی به چه میزان بود؟'
$ws.Range("A262").Value = 'همه موارد زیر در مورد پنوموسیت های نوع ii صحیح است، به جز:'
$ws.Range("A263").Value = 'شایع‌ترین علت خونریزی ثانویه پس از زایمان:-'
$ws.Range("A264").Value = 'شریان مهره‌ای از تمامی موارد زیر عبور می‌کند به جز:'
$ws.Range("A265").Value = 'کدام یک از موارد زیر از علل هیپوگلیسمی نیست؟'
$ws.Range("A266").Value = 'در طول جراحی کاشت حلزون چه چیزی قرار داده می‌شود؟'
$ws.Range("A267").Value = 'شاخص pma ارائه شده توسط شور و ماسلر روشی برای ارزیابی چیست؟'
$ws.Range("A268").Value = 'در عکس رادیوگرافی قفسه سینه، کدورت در ریه همراه با کلسیفیکاسیون نامنظم نشان دهنده کدام مورد است؟'
$ws.Range("A269").Value = 'کدام یک از موارد زیر باعث کاهش جذب استخوان در پوکی استخوان نمی‌شود؟'
$ws.Range("A270").Value = 'تومور دسموئید، درمان آن چیست؟'
$ws.Range("A271").Value = 'همه موارد زیر در مورد اشرشیاکلی انتروهموراژیک (ehec) صحیح هستند، به جز:'
$ws.Range("A272").Value = 'بر اساس دستگاه‌های تعیین موقعیت رأس کانال مقاومتی، مقاومتی که توسط رباط پریودنتال و غشای مخاطی دهان ارائه می‌شود چیست؟'
$ws.Range("A273").Value = 'همه موارد زیر در تعریف برونشیت مزمن صحیح هستند به جز:'
$ws.Range("A274").Value = 'گونیوسکوپی برای چه موردی استفاده می‌شود؟'
$ws.Range("A275").Value = 'در یک بیمار مبتلا به اسهال ناشی از ویبریو کلرا، کدام یک از موارد زیر مشاهده می‌شود؟  
الف) درد شکم  
ب) وجود گلبول‌های سفید در مدفوع  
ج) تب  
د) نوتروفیلی  
ه) وقوع موارد متعدد در یک منطقه'
$ws.Range("A276").Value = 'معیارهای دوک برای ارزیابی چه موردی استفاده می‌شود؟'
$ws.Range("A277").Value = 'آسکاریس لومبریکوئیدس باعث کمبود کدام ماده می‌شود؟'
$ws.Range("A278").Value = '10% تومور چیست؟'
$ws.Range("A279").Value = 'گاستن، یک پروتئین موجود در بزاق است'
$ws.Range("A280").Value = 'تمام موارد زیر در مورد دررفتگی مادرزادی لگن صحیح هستند به جز-'
$ws.Range("A281").Value = 'یک مرد 45 ساله پس از یک تصادف شدید رانندگی به بخش اورژانس منتقل می‌شود. معاینه فیزیکی نشان می‌دهد که بیمار دچار آسیب "straddle" به پرینه شده است. تصویربرداری mri نشان‌دهنده نشت ادرار و خون از مجرای بولبار پاره شده در شکاف سطحی پرینه است. کدام یک از فاشیاهای زیر مرزهای این فضا را تشکیل می‌دهند؟'
$ws.Range("A282").Value = 'کدام یک از موارد زیر ویژگی یک مهارکننده غیرانتخابی کاکس است اما ویژگی یک مهارکننده انتخابی کاکس 2 نیست؟'
$ws.Range("A283").Value = 'کدام گزینه در مورد بیاختیاری استرسی ادرار نادرست است؟'
$ws.Range("A284").Value = 'یک اختلال ژنتیکی باعث می‌شود که فروکتوز ۱،۶-بیفسفاتاز در کبد کمتر به تنظیم توسط فروکتوز ۲،۶-بیفسفات حساس باشد. تمام تغییرات متابولیک زیر در این اختلال مشاهده می‌شود به جز'
$ws.Range("A285").Value = 'رفلکس مچ پا توسط کدام ریشه عصبی منتقل می‌شود؟'
$ws.Range("A286").Value = 'سلول های t4 آنتی ژن را در ارتباط با تمام موارد زیر به جز کدام تشخیص می دهند؟ سپتامبر 2004'
$ws.Range("A287").Value = 'مری در 4 محل آناتومیک تنگ می‌شود. باریک‌ترین قسمت مری در کدام یک از تنگ‌شدگی‌های زیر قرار دارد؟'
$ws.Range("A288").Value = 'دوقلوهای مونوآمنیوتیک مونوکوریونیک در صورتی ایجاد می‌شوند که تقسیم پس از چه زمانی رخ دهد؟'
$ws.Range("A289").Value = 'شایع‌ترین علامت ایدز در یک نوزاد چیست؟'
$ws.Range("A290").Value = 'کدام یک از موارد زیر یک اختلال متابولیک اتوزومال غالب است؟'
$ws.Range("A291").Value = 'بر اساس قانون ثبت تولد و مرگ 1969، تولد و مرگ باید به ترتیب در چه مدت زمانی ثبت شوند؟'
$ws.Range("A292").Value = 'داروی انتخابی برای زنان باردار مبتلا به فشار خون مزمن که نیاز به درمان طولانی مدت با داروهای ضد فشار خون دارند: سپتامبر 2007'
$ws.Range("A293").Value = 'در کدام نوع از کارسینومای سلول کلیوی، هیستوپاتولوژی سلول‌های بزرگی با ظاهر گیاه‌مانند و هاله‌ی دور هسته‌ای دیده می‌شود؟'
$ws.Range("A294").Value = 'در مورد شریان ارتباطی خلفی کدام گزینه صحیح است؟'
$ws.Range("A295").Value = 'اندازه ذرات خطرناک که باعث پنوموکونیوز می‌شود از چه محدوده‌ای متغیر است؟'
$ws.Range("A296").Value = 'یک بیمار زن ۲۳ ساله با احساس ناتوانی در تنفس، سبکی سر و "گزگز" انگشتان دست، پا و اطراف دهان به اورژانس مراجعه می‌کند. این علائم هر بار که از تونل رانندگی می‌کند رخ می‌دهد و این بار نیز همین عامل محرک بوده است. کدام یک از مقادیر ph خون شریانی زیر با تشخیص او بیشترین تطابق را دارد؟'
$ws.Range("A297").Value = 'کدام تاندون در اندام فوقانی اغلب وجود ندارد؟'
$ws.Range("A298").Value = 'یک کودک 3 ماهه با استریدور متناوب مراجعه می‌کند. محتمل‌ترین علت چیست؟'
$ws.Range("A299").Value = 'لیگامان بری در تیروئید به کدام قسمت متصل می‌شود؟'
$ws.Range("A300").Value = 'اسپونجیوزیس شامل کدام لایه می‌شود؟'
$ws.Range("A301").Value = 'کاربرد دوره کمون همه موارد زیر است به جز؟'
$ws.Range("A302").Value = 'استنشاق پارچه آغشته به مواد مخدر چیست؟'
$ws.Range("A303").Value = 'نوکلئوزید از چه ترکیباتی تشکیل شده است؟  
1. پیریمیدین  
2. هیستون  
3. قند  
4. پورین  
5. فسفات'
$ws.Range("A304").Value = 'کدام یک از موارد زیر نشانه‌ای برای بیوپسی کلیه در یک کودک مبتلا به سندرم نفریتیک است؟'
$ws.Range("A305").Value = 'محرک اصلی برای ترشح سکرتین چیست؟'
$ws.Range("A306").Value = 'فردی پس از تجربه مشکل ناگهانی در تنفس به بخش اورژانس آورده می‌شود و یک تراکئوستومی پایین (زیر برآمدگی تیروئید) انجام می‌شود. کدام یک از عروق زیر ممکن است در این محل مشاهده شود؟'
$ws.Range("A307").Value = 'کمترین تأثیر بر انقباض پذیری میوکارد'
$ws.Range("A308").Value = 'تخم‌های بدون تقسیم در کدام یک از موارد زیر دیده می‌شوند؟'
$ws.Range("A309").Value = 'شایع ترین وضعیت پیش بدخیمی سرطان دهان کدام است؟'
$ws.Range("A310").Value = 'شایع‌ترین محل بروز نکروبیوزیس لیپوئیدیکا دیابتیکوروم کجاست؟'
$ws.Range("A311").Value = 'تمامی عبارات زیر در مورد تومورهای خارج مغزی صحیح هستند به جز -'
$ws.Range("A312").Value = 'بزرگ شدن لب که در عرض چند ثانیه تا 24 ساعت رخ می‌دهد، نشانه کدام مورد است؟'
$ws.Range("A313").Value = 'یک نوزاد مبتلا به زردی است که کف دست‌ها و پاها را نیز درگیر کرده است. سطح تقریبی بیلی‌روبین سرم چقدر است؟'
$ws.Range("A314").Value = 'کدام یک از موارد زیر یک اسید آمینه نیمه ضروری است؟'
$ws.Range("A315").Value = 'استریلیزاسیون فیبر نوری با چه روشی انجام می‌شود؟'
$ws.Range("A316").Value = 'بیمار با درد و تورم در واژن به شما مراجعه می‌کند. در معاینه، برآمدگی آبی رنگی در دیواره خلفی واژن پشت دهانه رحم مشاهده می‌شود. تاریخچه نشان می‌دهد که تورم با نزدیک شدن به قاعدگی بزرگ‌تر و دردناک‌تر می‌شود. تشخیص شما چیست؟'
$ws.Range("A317").Value = 'نظرسنجی ملی سلامت خانواده هر چند وقت یک‌بار انجام می‌شود؟'
$ws.Range("A318").Value = 'عبارت صحیح درباره نوروسیستیسرکوزیس چیست؟'
$ws.Range("A319").Value = 'دیپلوپی تک چشمی در همه موارد زیر رخ می دهد به جز:'
$ws.Range("A320").Value = 'سلول های کف آلود در کدام یک از موارد زیر دیده می شوند؟  
الف) سندرم آلپورت  
ب) بیماری نیمان پیک  
ج) آترواسکلروز  
د) پنومونی'
$ws.Range("A321").Value = 'در یک بیمار با آسیب سر، آسیب مغزی توسط کدام مورد تشدید می‌شود؟'
$ws.Range("A322").Value = 'کدام یک از گزینه‌های زیر در مورد نوار قلب در پریکاردیت حاد صحیح نیست؟'
$ws.Range("A323").Value = 'در چه سن حاملگی حداکثر افزایش برون ده قلبی رخ می‌دهد؟'
$ws.Range("A324").Value = 'کدام یک از موارد زیر با مصرف طولانی مدت باعث افزایش رنین می‌شود؟'
$ws.Range("A325").Value = 'کدام یک از موارد زیر برای میاستنی گراویس غیرمعمول است؟'
$ws.Range("A326").Value = 'یک زن خانه‌دار با پارونیشیای حاد از دو روز قبل مراجعه کرده است. بهترین عامل بی‌حسی برای پارونیشیای حاد کدام است؟'
$ws.Range("A327").Value = 'یک نوزاد پسر ترم ۵ روزه با تأخیر در دفع مکونیوم، اتساع شکم و استفراغ صفراوی مراجعه می‌کند. یافته‌های تنقیه باریم و بیوپسی روده در زیر نشان داده شده است. تشخیص چیست؟'
$ws.Range("A328").Value = 'آب مروارید هسته‌ای با کدام یک از موارد زیر مرتبط است؟'
$ws.Range("A329").Value = 'عضله ژنیوگلوسوس در سطح خلفی سیمفیزیس منتی به کدام ناحیه متصل می‌شود؟'
$ws.Range("A330").Value = 'پذیرفته‌شده‌ترین دلیل برای ازدحام دیررس دندان‌های پیشین چیست؟'
$ws.Range("A331").Value = 'به عنوان یک یافته منفرد، بالاترین خطر سندرم داون با کدام مورد مشاهده می‌شود؟'
$ws.Range("A332").Value = 'صلیب مالت از ویژگی‌های مشخصه کدام است؟'
$ws.Range("A333").Value = 'fnac (آسپیراسیون با سوزن نازک برای بررسی سلول‌شناسی) می‌تواند همه موارد زیر را تشخیص دهد به جز'
$ws.Range("A334").Value = 'کدام یک از شاخه‌های شریان فک بالا نیست؟'
$ws.Range("A335").Value = 'یکی از ویژگی‌های زیر می‌تواند برای تعریف لگن منقبض استفاده شود'
$ws.Range("A336").Value = 'همه موارد زیر در کم خونی فقر آهن کاهش می یابند به جز:'
$ws.Range("A337").Value = 'کم خونی همولیتیک با تست کومبس مثبت در همه موارد زیر دیده می‌شود به جز-'
$ws.Range("A338").Value = 'تومور مرتبط با بهترین پیش‌آگهی در کودکان کدام است؟'
$ws.Range("A339").Value = 'در کدام یک از بیماران زیر می‌توان اپیزیوتومی میانی را برتر از اپیزیوتومی میانی-جانبی در هنگام زایمان جنین در نظر گرفت؟'
$ws.Range("A340").Value = 'لارنگوسل از کجا نشأت می‌گیرد؟'
$ws.Range("A341").Value = 'کدام یک از داروهای زیر اخیراً برای درمان سرطان پروستات تأیید شده است؟ (تکرار)'
$ws.Range("A342").Value = 'برآورد سن از طریق معاینه دندان ها توسط کدام روش انجام می شود؟'
$ws.Range("A343").Value = 'همه موارد زیر در مورد سارکوم کاپوسی صحیح است به جز یک مورد:'
$ws.Range("A344").Value = 'غدد برونر در کدام قسمت دیده می‌شوند؟'
$ws.Range("A345").Value = 'برادی کاردی نسبی همراه با تب مشخصه کدام بیماری است؟'
$ws.Range("A346").Value = 'داروی انتخابی برای کارسینوم سلول کلیوی غیرقابل جراحی کدام است؟'
$ws.Range("A347").Value = 'گره راتر چیست؟'
$ws.Range("A348").Value = 'یک مرد ۶۵ ساله سیگاری با هماچوری کامل بدون درد مراجعه می‌کند. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A349").Value = 'کدام یک از اندام های زیر کمترین احتمال درگیری را در بیماری بهجت دارد؟'
$ws.Range("A350").Value = 'تمام موارد زیر در مورد آسپرین صحیح است به جز:'
$ws.Range("A351").Value = 'کدام یک از موارد زیر به طور مستقیم بر گیرنده‌های نیکوتینی استیل‌کولین (nm) اثر می‌گذارد؟'
$ws.Range("A352").Value = 'ژن hla hi کد کننده کدام مورد است؟'
$ws.Range("A353").Value = 'درگیری ماکولا در کدام یک از موارد زیر شایع است؟'
$ws.Range("A354").Value = 'سلول های حافظه t را می توان با استفاده از کدام نشانگر شناسایی کرد؟'
$ws.Range("A355").Value = 'پلی سیتمی ورا زمانی تشخیص داده می‌شود که هماتوکریت مطلق وریدی بیشتر از چه مقداری باشد؟'
$ws.Range("A356").Value = 'کدام گزینه برای درمان ریشه‌کنی پلاسمودیوم ویواکس استفاده می‌شود؟'
$ws.Range("A357").Value = 'کدام عضله دست دارای یک استخوان سزاموئید است؟'
$ws.Range("A358").Value = 'کدام یک از اعصاب زیر بیشتر عضلات ذاتی حنجره را عصب‌دهی می‌کند؟'
$ws.Range("A359").Value = 'خونریزی بینی تروماتیک در کدام مورد علت محسوب نمی‌شود؟'
$ws.Range("A360").Value = 'آسیب عصب مدین در مچ دست باعث همه موارد زیر می‌شود به جز -'
$ws.Range("A361").Value = 'سندرم لاوژیه-هانتزیکر با کدام مورد مرتبط است؟'
$ws.Range("A362").Value = 'کدام یک از عوامل زیر مهم‌ترین نقش را در پیشگیری از اندوفتالمیت در جراحی آب مروارید دارد؟'
$ws.Range("A363").Value = 'بیماری کرون معمولاً کدام قسمت را درگیر می‌کند؟'
$ws.Range("A364").Value = 'تخلیه لنفاوی قسمت خارجی فوقانی پستان عمدتاً به کدام ناحیه انجام می‌شود؟'
$ws.Range("A365").Value = 'تعداد رتیکولوسیت‌ها در زردی همولیتیک چقدر است؟'
$ws.Range("A366").Value = 'کولستئاتوما معمولاً کدام ساختار را سوراخ می‌کند؟'
$ws.Range("A367").Value = 'کدام یک از موارد زیر یک انتقال دهنده عصبی مونوآمین در مغز است؟'
$ws.Range("A368").Value = 'رتینوبلاستوما با کدام یک از موارد زیر مرتبط است؟'
$ws.Range("A369").Value = 'ضایعه سفید در حفره دهان در همه موارد زیر دیده می‌شود به جز؟'
$ws.Range("A370").Value = 'تیبیالیس خلفی به تمام استخوان های تارس متصل است، به جز:'
$ws.Range("A371").Value = 'تایپ فاژ به طور گسترده برای طبقه‌بندی درون‌گونه‌ای کدام یک از باکتری‌های زیر استفاده می‌شود؟'
$ws.Range("A372").Value = 'مسیر داخل عضلانی در هنگام مصرف کدام دارو منع شده است؟'
$ws.Range("A373").Value = 'مدیریت بیماران مبتلا به نوتروپنی تب‌دار تحت درمان ضدسرطان شامل همه موارد زیر به جز کدام است؟'
$ws.Range("A374").Value = 'پای راکر، میکروسفالی، پس سری برجسته و میکروگناتی از ویژگی‌های بالینی بارز کدام بیماری هستند؟'
$ws.Range("A375").Value = 'در مورد هیپوونتیلاسیون اولیه ریوی کدام گزینه صحیح است؟'
$ws.Range("A376").Value = 'فرم عفونی کننده برای پشه در پلاسمودیوم فالسیپاروم چیست؟'
$ws.Range("A377").Value = 'در یک مادر rh منفی که نوزاد rh مثبت به دنیا آورده است، پروپیلکتیک anti d در چه شرایطی تجویز می‌شود؟'
$ws.Range("A378").Value = 'همه موارد زیر گلیکوژنوز کبدی هستند به جز –'
$ws.Range("A379").Value = 'یک دختر ۴ ساله به مطب پزشک اطفال آورده می‌شود. پدرش گزارش می‌دهد که او ناگهان رنگ پریده شده و در حالی که او و سگ کوچکش را دنبال می‌کرده، از دویدن ایستاده است. پس از ۳۰ دقیقه، رنگ پریدگی او برطرف شده و می‌خواهد بازی را از سر بگیرد. او قبلاً هیچ‌گونه حمله‌ای نداشته و هرگز سیانوتیک نبوده است. معاینه فیزیکی او طبیعی بود و عکس قفسه سینه و اکوکاردیوگرام نیز طبیعی بودند. نوار قلب الگویی را نشان می‌دهد که در صفحه بعد مشاهده می‌شود و نشان‌دهنده کدام یک از موارد زیر است؟'
$ws.Range("A380").Value = 'بیمار با تب درجه پایین و کاهش وزن، کاهش حرکت قفسه سینه در سمت راست، کاهش فرامیتوس، کدری در پرکاشن و کاهش صدای تنفس در سمت راست دارد. نای به سمت چپ منحرف شده است. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A381").Value = 'یک مرد ۲۴ ساله در جریان یک رابطه جنسی با یک فاحشه به hiv مبتلا می‌شود. او بلافاصله پس از این رابطه هیچ علائمی مشاهده نمی‌کند. در طول مرحله نهفته بدون علامت عفونت، ویروس به طور فعال در حال تکثیر است و می‌توان آن را در ارتباط با کدام یک از موارد زیر یافت؟'
$ws.Range("A382").Value = 'بیماری مرتبط با مصرف بیش از حد ذرت -'
$ws.Range("A383").Value = 'سقوط روی پاشنه با شکستگی استخوان کالکانه معمولاً با کدام مورد همراه است؟'
$ws.Range("A384").Value = 'همه موارد زیر از ویژگی‌های ureaplasma urealyticum هستند به جز؟'
$ws.Range("A385").Value = '"علامت کوسمول" چیست؟'
$ws.Range("A386").Value = 'موارد صحیح درباره عفونت salmonella typhi در روده کدامند؟'
$ws.Range("A387").Value = 'کدام یک از موارد زیر یک فرآیند برگشت‌پذیر است؟'
$ws.Range("A388").Value = 'کدام یک از جمعیت‌های سلول‌های بنیادی زیر در مغز استخوان یافت نمی‌شود؟'
$ws.Range("A389").Value = 'برای درخواست جنون در دادگاه، ipc کدام است: ap 06؛ up 09؛ راجستان 11'
$ws.Range("A390").Value = 'یک پسر 15 ساله با کم خونی و زردی مراجعه کرده است. در معاینه، هموگلوبین او 6 گرم در دسی لیتر بود، سونوگرافی سنگ کیسه صفرا را نشان داد و اسمیر محیطی موارد زیر را نشان داد. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A391").Value = 'یک مرد ۶۲ ساله به مدت ۴ ماه کمردرد دارد. در معاینه فیزیکی هیچ یافته غیرطبیعی مشاهده نمی‌شود. آزمایش cbc نشان‌دهنده تعداد wbc 3700 در میکرولیتر، هموگلوبین ۱۰.۳ گرم در دسی‌لیتر، هماتوکریت ۳۱.۱٪، mcv 85 فمتولیتر و تعداد پلاکت ۱۱۰,۰۰۰ در میکرولیتر است. پروتئین توتال سرم او ۸.۵ گرم در دسی‌لیتر با آلبومین ۴.۱ گرم در دسی‌لیتر است. رادیوگرافی قفسه سینه هیچ ناهنجاری در قلب یا ریه‌ها نشان نمی‌دهد، اما چندین ناحیه روشن در مهره‌ها مشاهده می‌شود. آسپیرات مغز استخوان از استخوان جناغ انجام می‌شود و ماده‌ای ژله‌ای و قرمز تیره در سرنگ به دست می‌آید. کدام یک از انواع سلول‌های زیر به احتمال زیاد در بررسی میکروسکوپی این آسپیرات به تعداد زیاد دیده می‌شود؟'
$ws.Range("A392").Value = 'بیماری نشان داده شده در این تست عملکرد ریوی (pft) بهترین توصیف را با کدام گزینه دارد؟'
$ws.Range("A393").Value = 'کدام یک از موارد زیر پیش‌بینی‌کننده ناتوانی بلندمدت در بیماران مبتلا به ام‌اس نیست؟'
$ws.Range("A394").Value = 'شایع‌ترین علت داخل برونشی هموپتیزی چیست؟'
$ws.Range("A395").Value = 'در تکنیک نیروی متعادل مورد استفاده برای آماده‌سازی کانال:'
$ws.Range("A396").Value = 'روابط اپیپلوئیک فورامن نشان داده شده است. اپیپلوئیک فورامن چه چیزی را به هم متصل می‌کند؟'
$ws.Range("A397").Value = 'خراشیدگی‌های گسترده بر روی بدن یک عابر پیاده یافت می‌شود. علت آن چیست؟'
$ws.Range("A398").Value = 'آنوسمی رحمت‌آمیز در کدام مورد مشاهده می‌شود؟'
$ws.Range("A399").Value = 'در کدام یک از انواع فتق زیر، احشای شکمی بخشی از دیواره کیسه فتق را تشکیل می‌دهند؟'
$ws.Range("A400").Value = 'صفحه راهنما در پروتز پارسیل متحرک باید در کجا قرار گیرد؟'
$ws.Range("A401").Value = 'هذیان در کدام یک از موارد زیر دیده نمی‌شود؟'
$ws.Range("A402").Value = 'مدت زمان درمان جذام چند باسیلری چیست؟'
$ws.Range("A403").Value = 'مقدار کاهش کاسپ توصیه شده برای دریافت یک ترمیم فلز-سرامیک چقدر است؟'
$ws.Range("A404").Value = 'در پالپوتومی، داروی مورد استفاده برای حفظ پالپ کدام است؟'
$ws.Range("A405").Value = 'کدام اسید آمینه می‌تواند در ph خنثی پروتونه و دپروتونه شود؟'
$ws.Range("A406").Value = 'یک زن 35 ساله با تشخیص اختلال پرکاری تیروئید. فولیکول تیروئید او پوشیده شده از'
$ws.Range("A407").Value = 'حلقه شاتسکی در کجا قرار دارد؟'
$ws.Range("A408").Value = 'عامل القای انتخابی در آسم -'
$ws.Range("A409").Value = 'کیسه کوچک معده از چه چیزی محدود شده است؟'
$ws.Range("A410").Value = 'کدام یک از موارد زیر یک بیماری ویروسی نوظهور است؟'
$ws.Range("A411").Value = 'بیمار در کدام یک از شرایط زیر پس از عمل جراحی بدون علامت می‌شود؟'
$ws.Range("A412").Value = 'همه موارد زیر در مورد سندرم پاترسون کلی صحیح است به جز؟'
$ws.Range("A413").Value = 'کدام یک از شرایط زیر ممکن است نیاز به تراکئوستومی اورژانسی برای جلوگیری از خفگی داشته باشد؟'
$ws.Range("A414").Value = 'کدام یک از موارد زیر توهم کاذب (pseudohallucination) را از توهم (hallucination) متمایز می کند؟'
$ws.Range("A415").Value = 'اپیدیدیم توسط چه چیزی پوشیده شده است؟'
$ws.Range("A416").Value = 'مقاومت محیطی با کدام مورد رابطه معکوس دارد؟'
$ws.Range("A417").Value = 'درصد از دست رفتن سلول‌های اندوتلیال در طی جراحی اتوماتیک اندوتلیال کراتوپلاستی با جدا کردن دسمه (dsaek) چقدر است؟'
$ws.Range("A418").Value = 'عوارض جانبی لوزارتان شامل همه موارد زیر است به جز:'
$ws.Range("A419").Value = 'نور ایمن مورد استفاده در رادیوگرافی معمولی چه رنگی از نور را ساطع می‌کند؟'
$ws.Range("A420").Value = 'فولیکول‌های تیروئیدی فعال که در حال ترشح هستند، توسط کدام نوع اپیتلیوم پوشیده شده‌اند؟'
$ws.Range("A421").Value = 'اکتینومیست ها چیستند؟'
$ws.Range("A422").Value = '"انعطاف مومی شکل" ویژگی کدام یک از موارد زیر است؟'
$ws.Range("A423").Value = 'شایع‌ترین علت رینوره موکوپورولنت یک‌طرفه در کودک چیست؟'
$ws.Range("A424").Value = 'شایع‌ترین محل خونریزی مغزی کجاست؟'
$ws.Range("A425").Value = 'در زردی پساکبدی، غلظت بیلیروبین کونژوگه در خون بیشتر از بیلیروبین غیرکونژوگه است زیرا'
$ws.Range("A426").Value = 'بیماری جنگلی کیسانور (kfd) توسط چه چیزی منتقل می‌شود؟'
$ws.Range("A427").Value = 'دوز بارگیری لفلونوماید در آرتریت روماتوئید چقدر است؟'
$ws.Range("A428").Value = 'همه موارد زیر در سندرم نفروتیک کاهش می یابند به جز'
$ws.Range("A429").Value = 'لنفوم بورکیت کدام جابجایی را نشان می‌دهد؟'
$ws.Range("A430").Value = 'داروی انتخابی برای ضربان‌های نارس بطنی (vpc) ناشی از مسمومیت با دیژیتال چیست؟'
$ws.Range("A431").Value = 'در مورد کودک یک ساله با pda کدام گزینه صحیح است؟'
$ws.Range("A432").Value = 'آنتی‌بیوتیکی که اثر بلوک کننده عصبی-عضلانی را تقویت می‌کند-'
$ws.Range("A433").Value = 'استراتژی درمانی فعلی برای بیمارانی که به hiv مبتلا شده‌اند، یک رژیم چنددارویی به نام درمان ضد رتروویروسی بسیار فعال (haart) است. یکی از انواع داروهای مورد استفاده در این درمان، آنالوگ نوکلئوزید/نوکلئوتید است، مانند دیدانوزین. کدام یک از گزینه‌های زیر بهترین توصیف از مکانیسم عمل این داروها را ارائه می‌دهد؟'
$ws.Range("A434").Value = 'در میان آمیلوبلاستوما تک کیستیک، میزان عود بالا مربوط به کدام نوع است؟'
$ws.Range("A435").Value = 'تمام موارد زیر در مورد جراحی فست ترک صحیح است به جز:'
$ws.Range("A436").Value = 'در مورد تجویز آهن صحیح است'
$ws.Range("A437").Value = 'یک دختر 4 ساله همیشه حتی در فصل تابستان مجبور بود جوراب گرم بپوشد. در معاینه فیزیکی مشخص شد که او فشار خون بالا دارد و نبض فمورال او در مقایسه با نبض رادیال و کاروتید ضعیف است. رادیوگرافی قفسه سینه نشان دهنده فرورفتگی قابل توجه دنده‌ها در امتداد لبه‌های پایینی آن‌ها بود. این وضعیت به دلیل کدام مورد است؟'
$ws.Range("A438").Value = 'آنتاگونیست هپارین کدام است؟'
$ws.Range("A439").Value = 'اسهال عفونی توسط کدام یک از موارد زیر ایجاد می‌شود؟'
$ws.Range("A440").Value = 'آسیب به پونز تحتانی در حالی که پونز فوقانی سالم باقی می‌ماند، منجر به چه چیزی می‌شود؟'
$ws.Range("A441").Value = 'بخش شناور کوه یخ نشان‌دهنده چه چیزی است که پزشک در جامعه مشاهده می‌کند؟'
$ws.Range("A442").Value = 'روانپزشکی مشاوره-ارتباط (c-l) شامل چه چیزی است؟'
$ws.Range("A443").Value = 'سندرم پاترسون براون کلی با همه موارد زیر به جز کدامیک مشخص می‌شود؟'
$ws.Range("A444").Value = 'همه موارد زیر در مورد اندوکاردیت ترومبوتیک غیرباکتریایی صحیح است، به جز-'
$ws.Range("A445").Value = 'همه موارد زیر درباره asha درست هستند به جز:'
$ws.Range("A446").Value = 'تیوپنتون در کدام مورد منع مصرف دارد؟'
$ws.Range("A447").Value = 'زخم‌های تروفیک ناشی از کدام یک از موارد زیر هستند؟  
الف) جذام  
ب) بیماری بورگر  
ج) سیرنگومیلی  
د) ترومبوز ورید عمقی  
ه) واریس سیاهرگ‌ها'
$ws.Range("A448").Value = 'کدام هورمون‌ها در ناحیه گلومرولوزا دیده می‌شوند؟'
$ws.Range("A449").Value = 'نیروهای جویدنی ایجاد شده توسط پروتز کامل دندان تقریباً چقدر است؟'
$ws.Range("A450").Value = 'یک نوزاد 3 ساعته با آپنه در حال تهویه با کیسه و ماسک در 30 ثانیه گذشته است، اکنون تنفس خودبه خودی با ضربان قلب 110 در دقیقه نشان می‌دهد. مرحله بعدی باید باشد –'
$ws.Range("A451").Value = 'یک زن 40 ساله با سابقه سردرد و حالت تهوع ناگهانی که با استراحت و مسکن برطرف شده است، مراجعه می‌کند. بعداً برای چند روز تاری دید ایجاد شده است. در روز بستری، فلج عصب سوم همراه با سفتی گردن داشته است. دمای بدن 100 درجه فارنهایت بود. محتمل‌ترین تشخیص چیست؟'
$ws.Range("A452").Value = 'ظاهر "لوله سربی" در رادیوگرافی باریم انما تشخیص کدام بیماری است؟'
$ws.Range("A453").Value = 'عامل استنشاقی هپاتوتوکسیک'
$ws.Range("A454").Value = 'تشخیص آگاماگلوبولینمی وابسته به x در چه شرایطی باید مورد شک قرار گیرد؟'
$ws.Range("A455").Value = 'سندرم کوشینگ وابسته به غذا نتیجه بیان نابجای کدام مورد است؟'
$ws.Range("A456").Value = 'تمامی عبارات در مورد سندرم پلامر-وینسون صحیح هستند به جز -'
$ws.Range("A457").Value = 'کاهش تأخیر در شروع خواب rem ویژگی کدام یک از موارد زیر است؟'
$ws.Range("A458").Value = 'علامت "s طلایی" در کدام مورد دیده می‌شود؟'
$ws.Range("A459").Value = 'شایع‌ترین شکل هیستری تجزیه‌ای چیست؟'
$ws.Range("A460").Value = 'کدام یک از داروهای زیر برای درمان دیستونی گردن تأیید شده است؟'
$ws.Range("A461").Value = 'کدام توکسین استرپتوکوک باعث همولیز می‌شود؟'
$ws.Range("A462").Value = 'همه موارد زیر از یافته‌های رادیوگرافی در آرتریت روماتوئید هستند به جز:'
$ws.Range("A463").Value = 'سنگ های اگزالات در کدام مورد یافت می شوند؟'
$ws.Range("A464").Value = 'یک مرد بالغ متهم به تجاوز جنسی در مدت ۱۲ ساعت برای معاینه پزشکی آورده می‌شود. تجاوز جنسی در صورت وجود موارد زیر محرز نمی‌شود:'
$ws.Range("A465").Value = 'یک زن جوان دارای مقادیر آزمایشگاهی زیر است: هموگلوبین=9.8 گرم٪، mcv=70، آهن سرم=60، فریتین سرم=100، تشخیص چیست؟'
$ws.Range("A466").Value = 'افزایش سطح بتا-hcg در کدام یک از موارد زیر مشاهده می‌شود؟'
$ws.Range("A467").Value = 'کودکی با آدنوم سباسه و تشنج های مقاوم به درمان مراجعه کرده است. mri آستروسیتوم های سلول غول پیکر زیر اپاندیمال را نشان داد. علت چیست؟'
$ws.Range("A468").Value = 'd - sleep چیست؟'
$ws.Range("A469").Value = 'هنگام ارائه بیهوشی عمومی در دوران بارداری، حداقل غلظت آلوئولی (mac) چیست؟'
$ws.Range("A470").Value = 'کدام یک از موارد زیر باعث درد حاد مقعدی نمی‌شود؟'
$ws.Range("A471").Value = 'در مورد شالازیون با عودهای مکرر در همان محل، برش و کورتاژ باید:'
$ws.Range("A472").Value = 'نمایش ریاضی قانون مربع معکوس چیست؟'
$ws.Range("A473").Value = 'کدام استخوان بیشتر در بیماری پاژه درگیر می‌شود؟'
$ws.Range("A474").Value = 'در بارت کلاسیک، طول اپیتلیوم ستونی در مری دیستال چقدر است؟'
$ws.Range("A475").Value = 'ضایعه هدف یا عنبیه در کدام یک از موارد زیر دیده می‌شود؟'
$ws.Range("A476").Value = 'کودکی با بیرون زدگی حلقه‌های روده از یک نقص در دیواره شکم در سمت راست ناف مراجعه کرده است. تشخیص احتمالی چیست؟'
$ws.Range("A477").Value = 'کدام مورد در پلی سیتمی ورا دیده نمی‌شود؟'
$ws.Range("A478").Value = 'عوارض جانبی شبیه به سندرم اکستراپیرامیدال در کدام یک دیده می‌شود؟'
$ws.Range("A479").Value = 'چه کسری از عمر خود را با جویدن با دندان‌های شیری سپری می‌کنیم؟'
$ws.Range("A480").Value = 'منیزیم در کدام مورد نقش ندارد؟'
$ws.Range("A481").Value = 'برنامه ایندیرا آواس یوجانا تحت کدام یک از وزارتخانه‌های زیر است؟'
$ws.Range("A482").Value = 'یک بیمار مرد ۵۰ ساله از سنگینی صورت، سردرد و گرفتگی بینی شکایت دارد. بیمار سابقه سینوزیت مزمن دارد. عکس رادیوگرافی با نمای واترز گرفته شده است. وضعیت بیمار در این رادیوگرافی به چه صورت است؟'
$ws.Range("A483").Value = 'عرض یا عمق بیدینگ اتصال دهنده اصلی کام باید چقدر باشد؟'
$ws.Range("A484").Value = 'رنگ آمیزی دو قطبی مشخصه کدام یک از موارد زیر است؟'
$ws.Range("A485").Value = 'افزایش جریان خون ریوی که به صورت رادیولوژیکی با همه موارد زیر به جز کدامیک نشان داده می‌شود؟'
$ws.Range("A486").Value = 'همه موارد زیر به دررفتگی مکرر شانه مرتبط هستند به جز'
$ws.Range("A487").Value = 'معیارهای تشخیصی سندرم متابولیک شامل موارد زیر است (s)-'
$ws.Range("A488").Value = 'حلقه ویلیس در کجا قرار دارد؟'
$ws.Range("A489").Value = 'اسیدهای صفراوی از کجا بازجذب می‌شوند؟'
$ws.Range("A490").Value = 'سمپاتکتومی کمری برای همه موارد زیر به جز کدام انجام می‌شود؟ سپتامبر 2012'
$ws.Range("A491").Value = 'نقاط راث در کدام مورد دیده می‌شوند؟'
$ws.Range("A492").Value = 'قوی‌ترین و طولانی‌اثرترین عامل بیهوشی کدام است؟'
$ws.Range("A493").Value = 'اندیکاسیون ختنه در کودکان:'
$ws.Range("A494").Value = 'درگیری کلیوی در مولتیپل میلوما با کدام مورد مشخص می‌شود؟'
$ws.Range("A495").Value = 'تعداد مولکول‌های atp تولید شده توسط مسیر hmp چقدر است؟'
$ws.Range("A496").Value = 'تخلیه لنفاوی دیواره جانبی بینی'
$ws.Range("A497").Value = 'کدام گزینه در مورد واکسن sa-14-14-2 انسفالیت ژاپنی صحیح است؟'
$ws.Range("A498").Value = 'تمامی عبارات زیر در مورد آنتی‌ژن کربوهیدرات صحیح هستند، به جز:'
$ws.Range("A499").Value = 'کدام نیترات متابولیسم عبور اول را انجام نمی‌دهد؟'
$ws.Range("A500").Value = 'عبارت صحیح درباره ''نرخ باروری کل'' چیست؟'
$ws.Range("A501").Value = 'تمام گزینه‌های زیر در مورد فیلاریا صحیح هستند به جز کدام یک؟'
$ws.Range("A502").Value = 'سویه‌ای از e. coli که سم شیگا تولید می‌کند (ec o157:h7) را می‌توان با استفاده از کدام محیط کشت از سایر سویه‌های e. coli تشخیص داد؟'
$ws.Range("A503").Value = 'زخم مارژولین'
$ws.Range("A504").Value = 'نوزاد مرده به دنیا آمده به چه صورت تعریف می‌شود؟'
